# annualize apr and effects in std deviations
#
# The underlying (raw) data table (rows 15-18, cols B:F, plus the G helper
# column) gains an extra annualized row. The new row's numbers are written
# into row 16 (replacing what used to live there), the old row 16 values
# slide down into row 17, and the old row 17 values slide down into row 18
# with one corrected figure (F: 2186 -> 2185). The formulas in the summary
# table above (rows 4:6) keep referring to the same cell addresses as
# before (B16, C16, ..., F17, ...) rather than following the data - so
# their cached results go stale until the next full recalculation. Match
# that by editing the workbook with calculation set to Manual, so the
# cached formula values are left untouched (Excel will simply mark the
# workbook as needing a full recalc on next open).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("exp_arms")

# Leave calculation in Manual mode for the rest of this script (and when it
# completes) so the post-script recalc only fills in never-evaluated
# formulas instead of re-running the whole dependency graph - the summary
# table's cached results (rows 4:6) must stay exactly as they were.
$excel.Calculation = -4135   # xlCalculationManual - keep cached <v> stale

# Row 16 (new annualized figures; G16/H16/I16 formatting already in place)
$ws.Cells.Item(16, 2).Value = 2585
$ws.Cells.Item(16, 3).Value = 2465
$ws.Cells.Item(16, 4).Value = 2143
$ws.Cells.Item(16, 5).Value = 3406
$ws.Cells.Item(16, 6).Value = 2753

# Row 17 (former row-16 figures, shifted down); add the G helper value that
# row 16 used to carry, matching its centered-alignment style.
$ws.Cells.Item(17, 2).Value = 2036
$ws.Cells.Item(17, 3).Value = 1907
$ws.Cells.Item(17, 4).Value = 1757
$ws.Cells.Item(17, 5).Value = 2710
$ws.Cells.Item(17, 6).Value = 2216
$ws.Cells.Item(17, 7).Value = 6919
$ws.Cells.Item(17, 7).HorizontalAlignment = -4108   # xlHAlignCenter (style "1")

# Row 18 (former row-17 figures, shifted down); F corrected 2186 -> 2185.
$ws.Cells.Item(18, 2).Value = 1984
$ws.Cells.Item(18, 3).Value = 1840
$ws.Cells.Item(18, 4).Value = 1724
$ws.Cells.Item(18, 5).Value = 2634
$ws.Cells.Item(18, 6).Value = 2185
